$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used-range row extent so we cover every data row in column G
$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    # Find the (last) exact, case-sensitive occurrence of the token "System"
    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
        }
    }

    # Nothing to do if there's no exact "System" token, or it's already first
    if ($idx -le 0) {
        continue
    }

    $rest = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            $rest += $parts[$i]
        }
    }

    $newParts = @("System") + $rest
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
